# "winch switch is on cyrpus module yo."
#
# - Outputs sheet: rows 6/7 (winchMotor / winchSwitch) had "???" placeholders
#   in their Slot/Channel columns (C/D). Those get filled in with real slot 4
#   (and channels 4 / 3 respectively).
# - Other Inputs sheet: row 8 (winchSwitch) had "???" placeholders in C8/D8.
#   Those get replaced with a single centered, merged "On Cyprus" label.
# - Selections on both sheets move to reflect where the editor ended up.

$wb = $excel.ActiveWorkbook

# --- Outputs sheet --------------------------------------------------------
$wsOutputs = $wb.Worksheets.Item("Outputs")

$wsOutputs.Range("C6").Value = 4
$wsOutputs.Range("D6").Value = 4
$wsOutputs.Range("C7").Value = 4
$wsOutputs.Range("D7").Value = 3

[void]$wsOutputs.Range("A19").Select()

# --- Other Inputs sheet ----------------------------------------------------
$wsOther = $wb.Worksheets.Item("Other Inputs")

# Center-align first so the new shared style is created before the merge
# (matches the single new centered cellXfs entry added by the edit).
$wsOther.Range("C8:D8").HorizontalAlignment = -4108   # xlCenter
$wsOther.Range("C8").Value = "On Cyprus"
$wsOther.Range("D8").Value = ""
$wsOther.Range("C8:D8").Merge()

$wsOther.Range("C9").Select()
